# Auto-generated: apply cryptos.xlsx price/volume refresh per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.066.07"
$ws.Range("E2").Value = "  -5.35%  "
$ws.Range("D3").Value = "2.221.45"
$ws.Range("E3").Value = "  -6.65%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "323.47"
$ws.Range("E5").Value = "  -2.22%  "
$ws.Range("D6").Value = "98.74"
$ws.Range("E6").Value = "  -9.68%  "
$ws.Range("D7").Value = "0.581"
$ws.Range("E7").Value = "  -9.05%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "0.569"
$ws.Range("E9").Value = "  -8.07%  "
$ws.Range("D10").Value = "36.71"
$ws.Range("E10").Value = "  -11.27%  "
$ws.Range("E11").Value = "  -3.77%  "
$ws.Range("D12").Value = "0.0827"
$ws.Range("E12").Value = "  -10.28%  "
$ws.Range("D13").Value = "7.63"
$ws.Range("E13").Value = "  -10.81%  "
$ws.Range("E14").Value = "  -2.44%  "
$ws.Range("D15").Value = "2.561.06"
$ws.Range("E15").Value = "  -6.58%  "
$ws.Range("D16").Value = "0.859"
$ws.Range("E16").Value = "  -12.71%  "
$ws.Range("D17").Value = "14.32"
$ws.Range("E17").Value = "  -7.70%  "
$ws.Range("D18").Value = "2.226.66"
$ws.Range("D19").Value = "42.972.36"
$ws.Range("E19").Value = "  -5.44%  "
$ws.Range("D20").Value = "13.74"
$ws.Range("E20").Value = "  -9.97%  "
$ws.Range("D21").Value = "0.0₃0963"
$ws.Range("E21").Value = "  -9.83%  "
$ws.Range("D22").Value = "6.51"
$ws.Range("E22").Value = "  -11.39%  "
$ws.Range("D23").Value = "3.27"
$ws.Range("E23").Value = "  -11.62%  "
$ws.Range("D24").Value = "65.17"
$ws.Range("E24").Value = "  -11.32%  "
$ws.Range("D25").Value = "236.61"
$ws.Range("E25").Value = "  -10.60%  "
$ws.Range("D26").Value = "2.19"
$ws.Range("E26").Value = "  -5.72%  "
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("E28").Value = "  +1.46%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.00"
$ws.Range("E29").Value = "  -12.01%  "
$ws.Range("E30").Value = "  -3.23%  "
$ws.Range("D31").Value = "6.38"
$ws.Range("E31").Value = "  -14.68%  "
$ws.Range("D32").Value = "36.42"
$ws.Range("E32").Value = "  -2.08%  "
$ws.Range("D33").Value = "20.23"
$ws.Range("E33").Value = "  -10.15%  "
$ws.Range("D34").Value = "0.0863"
$ws.Range("E34").Value = "  -10.10%  "
$ws.Range("D35").Value = "155.37"
$ws.Range("E35").Value = "  -8.20%  "
$ws.Range("E36").Value = "  +0.65%  "
$ws.Range("D37").Value = "2.66"
$ws.Range("E37").Value = "  -7.07%  "
$ws.Range("D38").Value = "0.121"
$ws.Range("E38").Value = "  -8.52%  "
$ws.Range("D39").Value = "1.87"
$ws.Range("E39").Value = "  -5.30%  "
$ws.Range("D40").Value = "4.39"
$ws.Range("E40").Value = "  -7.86%  "
$ws.Range("E41").Value = "  -11.86%  "
$ws.Range("D42").Value = "3.68"
$ws.Range("E42").Value = "  -8.94%  "
$ws.Range("D43").Value = "0.0321"
$ws.Range("E43").Value = "  -10.06%  "
$ws.Range("D44").Value = "14.18"
$ws.Range("E44").Value = "  +8.83%  "
$ws.Range("E45").Value = "  -0.06%  "
$ws.Range("D46").Value = "1.728.07"
$ws.Range("E46").Value = "  -8.08%  "
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").Value = "0.202"
$ws.Range("E47").Value = "  -12.31%  "
$ws.Range("B48").Value = "BitcoinSV"
$ws.Range("C48").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D48").Value = "84.25"
$ws.Range("E48").Value = "  -14.18%  "
$ws.Range("D49").Value = "8.97"
$ws.Range("E49").Value = "  -4.25%  "
$ws.Range("D50").Value = "5.26"
$ws.Range("E50").Value = "  -14.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "74.20"
$ws.Range("E51").Value = "  -14.55%  "
